# chinh sua mot so gd
# Update "Loai san pham" (column E) to proper category names, and fix the
# garbled "Don vi tinh" (column G) value for food rows to "Phần".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column E (Loai san pham) - rows 4..29
$ws.Range("E4:E5").Value = "Bia"
$ws.Range("E6:E7").Value = "Nước ngọt"
$ws.Range("E8").Value = "Bia"
$ws.Range("E9:E10").Value = "Nước ngọt"
$ws.Range("E11:E24").Value = "Thức ăn"
$ws.Range("E25:E29").Value = "Đồ uống"

# Column G (Don vi tinh) - fix garbled "ph?n"/"Ph?n" -> "Phần"
$ws.Range("G11:G24").Value = "Phần"
